$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 13123
$ws.Range("B3").Value = 2412
$ws.Range("C3").Value = 241241

$ws.Range("C3").Select()
